$d = $word.ActiveDocument

$old0 = 'The merchant uses the app to scan the code and requests the buyer to enter a confirmation pin. After entering it on the merchant’s smartphone, the app connects to the blockchain to submit the transaction. If the code and pin match, and there are enough funds in the card address, the transaction is instantly verified and recorded on the blockchain. The process of scanning a card and entering a pin to confirm is already well familiar to current debit card users, and so it would be easy for everyone to get used to.'
$new0 = 'Торговец использует приложение для сканирования кода и просит покупателя ввести подтверждающий pin-код. После ввода кода покупателем на смартфоне продавца, приложение подключается к блокчейну для отправки транзакции. Если код и pin совпадают, и на адресе карты достаточно средств, транзакция мгновенно проверяется и записывается в блокчейне. Процесс сканирования карты и ввода PIN-кода для подтверждения уже хорошо знаком нынешним пользователями дебетовых карт, поэтому не вызовет затруднений.'
$found0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2)
Write-Host "Replace 0: $found0"

$old1 = 'A REVOLUTIONARY PAYMENT SYSTEM EVERYONE CAN LOVE'
$new1 = 'Революционная платёжная система для каждого'
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Replace 1: $found1"

$old2 = 'This ease of use and familiarity present the key advantage of this system over other crypto payments, which require both parties to be online for a transaction. Once the card is created and loaded, only the merchant is responsible for having a live connection to the internet and the consumer can take the funds anywhere without worrying about a dead battery or lack of signal in remote areas. No need for fumbling with mobile wallet apps, nor any memorization beyond the simple numerical pin just like the ones used for debit purchases. In addition to this, it also offers the advantage of being usable in places where traditional banking services are only partially available, or even absent entirely. Whereas other mobile payment solutions require bank accounts for both parties to send/receive the funds, this one uses the SmartCash blockchain instead and is thus accessible to anyone who wishes to use it.'
$new2 = 'Простота использования и очевидность процесса оплаты представляют собой ключевые преимущества этой системы над другими криптовалютными платежами; теперь нет необходимости, чтобы обе стороны были в сети для совершения транзакции. Как только карта будет создана и пополнена, только продавец несет ответственность за подключение к Интернету; потребитель же всегда будет иметь доступ к своим средствам — в любом месте, не беспокоясь о мертвой батарее или отсутствии сигнала в отдаленных районах. Не нужно возиться с мобильными приложениями для кошельков, не требуется запоминать коды или пароли; от сторон не потребуется ничего, что лежит за пределами использования обычных дебетовых карт. В дополнение к этому, SmartCard также предлагает преимущество использования в тех местах, где традиционные банковские услуги доступны только частично или вообще отсутствуют. В то время как другие мобильные платежные решения требуют наличие банковских счетов обеих сторон для отправки / получения средств, SmartCash использует собственный блокчейн, предлагая мгновенные и удобные платежи абсолютно для всех.'
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Replace 2: $found2"

$old3 = 'Naturally, these benefits would also be quite attractive to buyers and sellers in more developed nations as well. Other crypto solutions fall somewhat short of their promises. Crypto debit cards suffer from the same fees and transaction delays as any other current plastic card, and sending money through current mobile wallets requires a lot of time and effort to make the payments. Compared to these options, SmartCard payments are in a league of their own in convenience and speed. Add to this the fact that the SmartCash community is actively seeking, recruiting, and sponsoring business owners and entrepreneurs globally, and it is easy to imagine how the SmartCard could increasingly gain global adoption across a wide range of markets and economies.'
$new3 = 'Конечно, эти особенности также будут весьма привлекательными для предпринимателей и потребителей в более развитых странах. Другие крипто-решения несколько отстают от своих обещаний. “Криптовалютные” дебетовые карты страдают от таких же сборов и задержек транзакций, как и любые другие пластиковые банковские карты, а отправка денег через актуальные мобильные кошельки требует немало времени и усилий для совершения даже самых простых платежей. По сравнению с этим, платежи через SmartCard находятся в высшей лиге по удобству и скорости. Также не стоит забывать, что сообщество SmartCash активно ищет, привлекает и спонсирует владельцев бизнеса и предпринимателей во всем мире, что, в недалеком будущем, значительно упростит внедрение SmartCard в экономические отношения на многих рынках.'
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Host "Replace 3: $found3"

$old4 = '5 REASONS TO CHOOSE SMARTCARD'
$new4 = '5 причин выбрать SmartCard'
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Host "Replace 4: $found4"

$old5 = 'NO FEE'
$new5 = 'Отсутствие комиссий'
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Host "Replace 5: $found5"

$old6 = '0% transaction fees with blockchain fee of less than 1/10 of a cent.'
$new6 = '0% - комиссия за транзакцию; менее 1/10 цента - комиссия сети блокчейн.'
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
Write-Host "Replace 6: $found6"

$old7 = 'SECURE'
$new7 = 'Безопасность'
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
Write-Host "Replace 7: $found7"

$old8 = 'Manage your money in a way that is completely secure'
$new8 = 'Просто используйте карту. Ваши средства надежно защищены.'
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
Write-Host "Replace 8: $found8"

$old9 = 'CARDS'
$new9 = 'Карты'
$found9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
Write-Host "Replace 9: $found9"

$old10 = 'Easy self made cards anyone can create'
$new10 = 'Удобные карты, которые может создать каждый.'
$found10 = $d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, 1, $false, $new10, 2)
Write-Host "Replace 10: $found10"

$old11 = 'CONVENIENCE'
$new11 = 'Удобство'
$found11 = $d.Content.Find.Execute($old11, $true, $false, $false, $false, $false, $true, 1, $false, $new11, 2)
Write-Host "Replace 11: $found11"

$old12 = 'Take funds anywhere without worrying about an internet connection'
$new12 = 'Где бы вы ни были, ваши средства всегда с вами.'
$found12 = $d.Content.Find.Execute($old12, $true, $false, $false, $false, $false, $true, 1, $false, $new12, 2)
Write-Host "Replace 12: $found12"

$old13 = 'SPEED'
$new13 = 'Скорость'
$found13 = $d.Content.Find.Execute($old13, $true, $false, $false, $false, $false, $true, 1, $false, $new13, 2)
Write-Host "Replace 13: $found13"

$old14 = 'Confirmation speed measured in fractions of a second'
$new14 = 'Подтверждение платежа происходит за доли секунды.'
$found14 = $d.Content.Find.Execute($old14, $true, $false, $false, $false, $false, $true, 1, $false, $new14, 2)
Write-Host "Replace 14: $found14"

$old15 = 'THE SMARTCARD INTRODUCTION'
$new15 = 'SmartCard: Введение'
$found15 = $d.Content.Find.Execute($old15, $true, $false, $false, $false, $false, $true, 1, $false, $new15, 2)
Write-Host "Replace 15: $found15"

$old16 = 'We’ve prepared an short video to help explain SmartCard in 3 minutes. Share with your friends!'
$new16 = 'Мы подготовили короткое видео, чтобы объяснить, что такое SmartCard. Поделитесь!'
$found16 = $d.Content.Find.Execute($old16, $true, $false, $false, $false, $false, $true, 1, $false, $new16, 2)
Write-Host "Replace 16: $found16"

$allFound = @($found0, $found1, $found2, $found3, $found4, $found5, $found6, $found7, $found8, $found9, $found10, $found11, $found12, $found13, $found14, $found15, $found16)
$successCount = ($allFound | Where-Object { $_ }).Count
Write-Host "Total replacements: $successCount / $($allFound.Count)"
